$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update position values (x1, y1, x2, y2) in rows 2 and 3
$ws.Range("A2").Value = -0.5
$ws.Range("B2").Value = -0.5
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0

$ws.Range("A3").Value = 0.5
$ws.Range("B3").Value = 0
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = -0.5

# Update active cell selection to F11
$ws.Range("F11").Select()
